$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the selection for the bottomLeft pane
$ws.Range("X3").Select()

# Apply header row (row 1) style to the newly-used AX1 cell before writing its value
$ws.Range("AW1").Copy()
$ws.Range("AX1").PasteSpecial(-4122)

# Apply data row (row 2) style to the newly-used AX2 cell before writing its value
$ws.Range("AW2").Copy()
$ws.Range("AX2").PasteSpecial(-4122)

# Row 1 header values (Z1:AX1) - a new "inventory_item_id" column was inserted at Z,
# shifting metafield_* headers one column to the right
$ws.Range("Z1").Value = "inventory_item_id"
$ws.Range("AA1").Value = "metafield_Brand_id"
$ws.Range("AB1").Value = "metafield_Reference_Number_id"
$ws.Range("AC1").Value = "metafield_Collection_id"
$ws.Range("AD1").Value = "metafield_Product_Name_id"
$ws.Range("AE1").Value = "metafield_Description_id"
$ws.Range("AF1").Value = "metafield_Price_id"
$ws.Range("AG1").Value = "metafield_Price_Note_id"
$ws.Range("AH1").Value = "metafield_Quantity_id"
$ws.Range("AI1").Value = "metafield_Gender_id"
$ws.Range("AJ1").Value = "metafield_Function_id"
$ws.Range("AK1").Value = "metafield_Category_id"
$ws.Range("AL1").Value = "metafield_Sub_Category_id"
$ws.Range("AM1").Value = "metafield_Metal_id"
$ws.Range("AN1").Value = "metafield_Stone_id"
$ws.Range("AO1").Value = "metafield_Diamond_Clarity_id"
$ws.Range("AP1").Value = "metafield_Total_Diamond_Weight_id"
$ws.Range("AQ1").Value = "metafield_Also_Available_In_id"
$ws.Range("AR1").Value = "metafield_Width_id"
$ws.Range("AS1").Value = "metafield_Head_Size_id"
$ws.Range("AT1").Value = "metafield_Pearl_Type_id"
$ws.Range("AU1").Value = "metafield_Pearl_Size_id"
$ws.Range("AV1").Value = "metafield_Picture_1_id"
$ws.Range("AW1").Value = "metafield_Picture_2_id"
$ws.Range("AX1").Value = "metafield_Picture_3_id"

# Row 2 data values (Y2:AX2) - updated product/inventory/metafield ids
$ws.Range("Y2").Value = "7492350705833"
$ws.Range("Z2").Value = "44332328911017"
$ws.Range("AA2").Value = "20897912357033"
$ws.Range("AB2").Value = "20897912389801"
$ws.Range("AC2").Value = "20897912422569"
$ws.Range("AD2").Value = "20897912455337"
$ws.Range("AE2").Value = "20897912488105"
$ws.Range("AF2").Value = "20897912520873"
$ws.Range("AG2").Value = "20897912553641"
$ws.Range("AH2").Value = "20897912586409"
$ws.Range("AI2").Value = "20897912619177"
$ws.Range("AJ2").Value = "20897912651945"
$ws.Range("AK2").Value = "20897912684713"
$ws.Range("AL2").Value = "20897912717481"
$ws.Range("AM2").Value = "20897912783017"
$ws.Range("AN2").Value = "20897912815785"
$ws.Range("AO2").Value = "20897912848553"
$ws.Range("AP2").Value = "20897912881321"
$ws.Range("AQ2").Value = "20897912914089"
$ws.Range("AR2").Value = "20897912946857"
$ws.Range("AS2").Value = "20897912979625"
$ws.Range("AT2").Value = "20897913012393"
$ws.Range("AU2").Value = "20897913045161"
$ws.Range("AV2").Value = "20897913077929"
$ws.Range("AW2").Value = "20897913110697"
$ws.Range("AX2").Value = "20897913143465"
